# Auto-generated Excel COM-interop script applying scheduled-runner price updates
# to the Siren_Profits leve-profit tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 529.4286
$ws.Range("I2").Value = 651.2222
$ws.Range("J2").Value = 310.2
$ws.Range("K2").Value = 651.2222
$ws.Range("L2").Value = 310.2
$ws.Range("M2").Value = -538.2222
$ws.Range("N2").Value = -536.2
$ws.Range("H9").Value = 154.83333
$ws.Range("I9").Value = 200.5
$ws.Range("J9").Value = 63.5
$ws.Range("K9").Value = 200.5
$ws.Range("L9").Value = 63.5
$ws.Range("M9").Value = -31.5
$ws.Range("N9").Value = -401.5
$ws.Range("H19").Value = 1208.6316
$ws.Range("J19").Value = 1465
$ws.Range("L19").Value = 1465
$ws.Range("N19").Value = -1815
$ws.Range("H74").Value = 4697.3335
$ws.Range("I74").Value = 4347.875
$ws.Range("K74").Value = 4347.875
$ws.Range("M74").Value = -3411.875
$ws.Range("H77").Value = 4697.3335
$ws.Range("I77").Value = 4347.875
$ws.Range("K77").Value = 21739.375
$ws.Range("M77").Value = -17059.375
$ws.Range("H98").Value = 69424.3
$ws.Range("I98").Value = 110794.6
$ws.Range("J98").Value = 28054
$ws.Range("K98").Value = 110794.6
$ws.Range("L98").Value = 28054
$ws.Range("M98").Value = -109296.6
$ws.Range("N98").Value = -31050
$ws.Range("H100").Value = 6099228
$ws.Range("I100").Value = 8561.200000000001
$ws.Range("J100").Value = 18280562
$ws.Range("K100").Value = 8561.200000000001
$ws.Range("L100").Value = 18280562
$ws.Range("M100").Value = -8020.200000000001
$ws.Range("N100").Value = -18281644
$ws.Range("H122").Value = 69424.3
$ws.Range("I122").Value = 110794.6
$ws.Range("J122").Value = 28054
$ws.Range("K122").Value = 332383.8
$ws.Range("L122").Value = 84162
$ws.Range("M122").Value = -329933.8
$ws.Range("N122").Value = -89062
$ws.Range("H125").Value = 4373.75
$ws.Range("I125").Value = 4500
$ws.Range("J125").Value = 4331.6665
$ws.Range("K125").Value = 40500
$ws.Range("L125").Value = 38984.9985
$ws.Range("M125").Value = -38040
$ws.Range("N125").Value = -43904.9985
$ws.Range("H132").Value = 3678.6667
$ws.Range("I132").Value = 3603.138
$ws.Range("J132").Value = 4226.25
$ws.Range("K132").Value = 10809.414
$ws.Range("L132").Value = 12678.75
$ws.Range("M132").Value = -8279.414000000001
$ws.Range("N132").Value = -17738.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5351.05
$ws.Range("I32").Value = 5448.4736
$ws.Range("J32").Value = 3500
$ws.Range("K32").Value = 5448.4736
$ws.Range("L32").Value = 3500
$ws.Range("M32").Value = -5161.4736
$ws.Range("N32").Value = -4074
$ws.Range("H45").Value = 12987.6875
$ws.Range("I45").Value = 16499
$ws.Range("J45").Value = 7135.5
$ws.Range("K45").Value = 16499
$ws.Range("L45").Value = 7135.5
$ws.Range("M45").Value = -16122
$ws.Range("N45").Value = -7889.5
$ws.Range("H46").Value = 10209.25
$ws.Range("J46").Value = 10953.714
$ws.Range("L46").Value = 10953.714
$ws.Range("N46").Value = -11591.714
$ws.Range("H110").Value = 1830.25
$ws.Range("I110").Value = 1830.25
$ws.Range("K110").Value = 1830.25
$ws.Range("M110").Value = 214.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2728.2917
$ws.Range("I94").Value = 2161.5
$ws.Range("J94").Value = 3861.875
$ws.Range("K94").Value = 2161.5
$ws.Range("L94").Value = 3861.875
$ws.Range("M94").Value = -1710.5
$ws.Range("N94").Value = -4763.875
$ws.Range("H99").Value = 25148.092
$ws.Range("I99").Value = 30292.111
$ws.Range("K99").Value = 30292.111
$ws.Range("M99").Value = -28794.111
$ws.Range("H107").Value = 4037.2104
$ws.Range("I107").Value = 3862.1538
$ws.Range("K107").Value = 3862.1538
$ws.Range("M107").Value = -1942.1538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H16").Value = 2171.2
$ws.Range("J16").Value = 2331
$ws.Range("L16").Value = 2331
$ws.Range("N16").Value = -2905
$ws.Range("H31").Value = 3406.561
$ws.Range("J31").Value = 4463.1816
$ws.Range("L31").Value = 4463.1816
$ws.Range("N31").Value = -5053.1816
$ws.Range("H34").Value = 3406.561
$ws.Range("J34").Value = 4463.1816
$ws.Range("L34").Value = 4463.1816
$ws.Range("N34").Value = -4867.1816
$ws.Range("H94").Value = 3049.5
$ws.Range("J94").Value = 1613
$ws.Range("L94").Value = 1613
$ws.Range("N94").Value = -2515
$ws.Range("H107").Value = 62514980
$ws.Range("I107").Value = 90929896
$ws.Range("K107").Value = 90929896
$ws.Range("M107").Value = -90927976
$ws.Range("H113").Value = 2171.2
$ws.Range("J113").Value = 2331
$ws.Range("L113").Value = 2331
$ws.Range("N113").Value = -6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 2087.5833
$ws.Range("J46").Value = 4312.2
$ws.Range("L46").Value = 12936.6
$ws.Range("N46").Value = -13118.6
$ws.Range("H56").Value = 5867.273
$ws.Range("I56").Value = 5867.273
$ws.Range("K56").Value = 5867.273
$ws.Range("M56").Value = -5337.273
$ws.Range("H122").Value = 1377.3529
$ws.Range("J122").Value = 1736.6666
$ws.Range("L122").Value = 15629.9994
$ws.Range("N122").Value = -20529.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 29886.75
$ws.Range("J52").Value = 29886.75
$ws.Range("L52").Value = 29886.75
$ws.Range("N52").Value = -30404.75
$ws.Range("H97").Value = 8755.037
$ws.Range("J97").Value = 2376.8572
$ws.Range("L97").Value = 2376.8572
$ws.Range("N97").Value = -3368.8572
$ws.Range("H132").Value = 3392.8262
$ws.Range("I132").Value = 2377.375
$ws.Range("J132").Value = 5713.857
$ws.Range("K132").Value = 7132.125
$ws.Range("L132").Value = 17141.571
$ws.Range("M132").Value = -4602.125
$ws.Range("N132").Value = -22201.571
$ws.Range("H134").Value = 41231.6
$ws.Range("J134").Value = 41231.6
$ws.Range("L134").Value = 123694.8
$ws.Range("N134").Value = -128764.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 26808.2
$ws.Range("I45").Value = 21347
$ws.Range("K45").Value = 21347
$ws.Range("M45").Value = -20940
$ws.Range("H61").Value = 4666.3335
$ws.Range("J61").Value = 4666.3335
$ws.Range("L61").Value = 4666.3335
$ws.Range("N61").Value = -5070.3335
$ws.Range("H113").Value = 4666.3335
$ws.Range("J113").Value = 4666.3335
$ws.Range("L113").Value = 4666.3335
$ws.Range("N113").Value = -9006.333500000001
$ws.Range("H122").Value = 3715.9583
$ws.Range("I122").Value = 3163.182
$ws.Range("J122").Value = 4183.6924
$ws.Range("K122").Value = 9489.545999999998
$ws.Range("L122").Value = 12551.0772
$ws.Range("M122").Value = -7039.545999999998
$ws.Range("N122").Value = -17451.0772
$ws.Range("H132").Value = 1412647.1
$ws.Range("I132").Value = 1966073
$ws.Range("K132").Value = 5898219
$ws.Range("M132").Value = -5895689

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 79991.664
$ws.Range("J46").Value = 79991.664
$ws.Range("L46").Value = 79991.664
$ws.Range("N46").Value = -80453.664
$ws.Range("H107").Value = 19908.883
$ws.Range("I107").Value = 1799.7693
$ws.Range("K107").Value = 5399.3079
$ws.Range("M107").Value = -3479.3079
$ws.Range("H113").Value = 8445.909
$ws.Range("I113").Value = 6966.6665
$ws.Range("K113").Value = 20899.9995
$ws.Range("M113").Value = -18729.9995
$ws.Range("H122").Value = 23920.27
$ws.Range("I122").Value = 2346.35
$ws.Range("K122").Value = 7039.049999999999
$ws.Range("M122").Value = -4589.049999999999
$ws.Range("H132").Value = 15435.5
$ws.Range("I132").Value = 25161.062
$ws.Range("K132").Value = 75483.186
$ws.Range("M132").Value = -72953.186
$ws.Range("H133").Value = 79997.5
$ws.Range("J133").Value = 79997.5
$ws.Range("L133").Value = 79997.5
$ws.Range("N133").Value = -90117.5
$ws.Range("H134").Value = 79991.664
$ws.Range("J134").Value = 79991.664
$ws.Range("L134").Value = 239974.992
$ws.Range("N134").Value = -245044.992
